$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = $null
$ws.Range("D3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("E3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H3").Value = $null
$ws.Range("I3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K3").Value = $null
$ws.Range("C4").Value = $null
$ws.Range("D4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("E4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H4").Value = $null
$ws.Range("I4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K4").Value = $null
$ws.Range("C5").Value = $null
$ws.Range("D5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("E5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H5").Value = $null
$ws.Range("I5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K5").Value = $null
$ws.Range("C6").Value = $null
$ws.Range("D6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("E6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("H6").Value = "8.300,01 TL - 199,41 TL"
$ws.Range("I6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("K6").Value = $null
$ws.Range("C8").Value = $null
$ws.Range("D8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("E8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("H8").Value = $null
$ws.Range("I8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("K8").Value = $null
$ws.Range("C9").Value = $null
$ws.Range("D9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("E9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("H9").Value = $null
$ws.Range("I9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("K9").Value = $null
$ws.Range("C10").Value = $null
$ws.Range("D10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("E10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("H10").Value = $null
$ws.Range("I10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("K10").Value = $null
$ws.Range("C11").Value = $null
$ws.Range("D11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("E11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("H11").Value = $null
$ws.Range("I11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("K11").Value = $null
$ws.Range("C12").Value = $null
$ws.Range("D12").Value = "WU: 0,75 USD–12 USD; Diğer: 700 TL–4.000 TL"
$ws.Range("K12").Value = $null
$ws.Range("C13").Value = $null
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 1.114 TL"
$ws.Range("H13").Value = "Hesaba: Asgari 1 TL | Azami 7,97 TL"
$ws.Range("I13").Value = "Hesaba: Asgari 1 TL | Azami 6,09 TL"
$ws.Range("K13").Value = $null
$ws.Range("C14").Value = $null
$ws.Range("E14").Value = "1.660 TL - 1.660 TL"
$ws.Range("H14").Value = "3.000 TL - 6.000 TL"
$ws.Range("K14").Value = $null
$ws.Range("K15").Value = "%0,3 Asgari Tutar: 237,26 TL Azami Tutar: 237,26 TL / 298,96 TL"
$ws.Range("H17").Value = " Asgari Tutar: 210 TL Azami Tutar: 210 TL / 210 TL / 210 TL"
$ws.Range("K17").Value = "%0,6 Asgari Tutar: 237,26 TL Azami Tutar: 237,26 TL / 3.034,67 TL"
$ws.Range("H20").Value = "100 TL"
$ws.Range("K20").Value = "147,11 TL"
$ws.Range("H21").Value = "%0,5 Asgari Tutar: 2.750 TL Azami Tutar: 2.750 TL"
$ws.Range("K21").Value = "%0,9 Asgari Tutar: 446,06 TL Azami Tutar: 446,06 TL / 2.427,26 TL"
$ws.Range("H22").Value = "%0,7 Asgari Tutar: 400 TL Azami Tutar: 400 TL / 400 TL"
$ws.Range("K22").Value = "%0,3 Asgari Tutar: 73,56 TL Azami Tutar: 73,56 TL / 9.115,86 TL"
$ws.Range("H23").Value = "57,5 TL"
$ws.Range("K23").Value = "64,8 TL"
$ws.Range("H24").Value = "350 TL"
$ws.Range("K24").Value = "446,06 TL"
$ws.Range("H25").Value = "400 TL"
$ws.Range("K25").Value = "374,4 TL"
